$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reconciliation rows (306-322) to append at the bottom of the sheet
$data = @(
    @(237676840777, "IVANS FANWOUM NOUPOUEH", "Rte_0", "Pk8", 471586.4666666666, 193917, -277669.4666666666, 0.411201367525814, "Ndogbong"),
    @(237677833877, "MEDJOM TAGNE MICHELLE GUILENE LA NEGRESSE SARL", "Rte_4", "Pk8", 38396.7, 131184, 92787.3, 3.416543609216417, "Ndogbong"),
    @(237678854978, "NSAMO NDJOUOHOU MICRANGE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_5", "Pk8", 40520, 92731, 52211, 2.288524185587364, "Ndogbong"),
    @(237679422591, "ETS LE CONTENT 42", "Rte_0", "Pk8", 141149.8, 3900, -137249.8, 0.02763021980902559, "Ndogbong"),
    @(237650353920, "MENIAPI HELENE EDOSSINE TOP MOBIL TELECOM", "Rte_4", "Socaver Ndongbong", 218865, 612128, 393263, 2.796829095561191, "Ndogbong"),
    @(237651927448, "ODETTE MABAKOU EPOUSE KENNE", "Rte_4", "Socaver Ndongbong", 128653.3333333333, 172411, 43757.66666666669, 1.340120737900301, "Ndogbong"),
    @(237653294562, "NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS", "Rte_5", "Socaver Ndongbong", 126220, 847745, 721525, 6.716407859293297, "Ndogbong"),
    @(237678046498, "MFS SOCAVER", "Rte_4", "Socaver Ndongbong", 41943.17, 6477, -35466.17, 0.1544232350582944, "Ndogbong"),
    @(237679428698, "ETS LE CONTENT 29", "Rte_8", "Socaver Ndongbong", 161427.5, 42, -161385.5, 0.0002601787180003407, "Ndogbong"),
    @(237679551262, "LA NEGRESSE LTDLA CBOX R1 MEGNE JUDITH", "Rte_5", "Socaver Ndongbong", 65595, 152434, 86839, 2.323866148334477, "Ndogbong"),
    @(237680574202, "TOUMEWO SAMUEL", "Rte_6", "Socaver Ndongbong", 86620, 224117, 137497, 2.587358577695682, "Ndogbong"),
    @(237681118330, "SAHA NDESA JONAS LTDLA_POLAS_OTH_NDOGBONG SERIE", "Rte_5", "Socaver Ndongbong", 152199.5714285714, 132896, -19303.57142857142, 0.8731693443852386, "Ndogbong"),
    @(237674446293, "ARSENE TITCHO KWAKEP", "Rte_0", "Total Ndokotti", 5994.285714285715, 11136, 5141.714285714285, 1.857769304099142, "Ndogbong"),
    @(237679085953, "FERNANDEZ NJOFANG TCHIYADJE", "Rte_0", "Total Ndokotti", 8950, 58169, 49219, 6.499329608938547, "Ndogbong"),
    @(237681662761, "ROLCHILE DJAMEN KOUDJOU", "Rte_0", "Total Ndokotti", 29783.82, 8697, -21086.82, 0.2920041821364755, "Ndogbong"),
    @(237682975726, "LA NEGRESSE SARL NYOUNG JOSEPH CLOTAIRE", "Rte_0", "Total Ndokotti", 29809.16666666666, 114072, 84262.83333333334, 3.826742333174919, "Ndogbong"),
    @(237683075075, "GLADYS LANG NGOINSEH", "Rte_0", "Total Ndokotti", 173418.2, 134739, -38679.19999999998, 0.7769599730593445, "Ndogbong"),
)

$startRow = 306
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
